$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Litr1-3"
$ws.Range("C3").Value = "Phys1-3"

$ws.Range("A6").Value = "English2-1"
$ws.Range("C6").Value = "Math2-3"
$ws.Range("E6").Value = "Phys2-5"
$ws.Range("A7").Value = "Math2-1"
$ws.Range("C7").Value = "Phys2-3"
$ws.Range("E7").Value = "Math2-5"

$ws.Range("B10").Value = "Russian3-2"
$ws.Range("C10").Value = "English3-3"
$ws.Range("D10").Value = "Litra3-4"
$ws.Range("E10").Value = "Phys3-5"
$ws.Range("B11").Value = "Phys3-2"
$ws.Range("D11").Value = "Phys3-4"

$ws.Range("C15").Value = "Phys4-3"
$ws.Range("C16").Value = "Math4-3"

$ws.Range("E19").Value = "Math5-5"
